$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the example data row (old row 5) and the Daily Totals row (old row 6)
$ws.Range("A5:M6").EntireRow.Delete()

$ws.Range("C11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("E13").Value = ""
$ws.Range("H13").Value = " ---------- 0.00"
$ws.Range("J13").Value = ""
$ws.Range("L13").Value = ""

$ws.Range("A2").Value = "Franchisee: Amto.Robert (Arthur Murray Thousand Oaks)"
$ws.Range("H2").Value = "(10/19/2025 - 10/25/2025)"
$ws.Range("K2").Value = "Week # 42"

$ws.Range("L13").Select()

Write-Host "done"
